# Generate Report for Handoff
# Re-run of the status report: the row that was "faec0f94-5a0c-4cd5-a6b8-19c4b111ec72.md"
# (previously "Handed back: in sync with en-US") moves to the bottom of each
# table and flips to "Ready for handoff" with fresh timestamps; the other two
# rows shift up to take its place.

$wb = $excel.ActiveWorkbook

function Set-CellAndLink {
    param(
        $ws,
        [int]$Row,
        [int]$Col,
        [string]$Value
    )
    $cell = $ws.Cells.Item($Row, $Col)
    $cell.Value = $Value
    foreach ($h in $ws.Hyperlinks) {
        if ($h.Range.Row -eq $Row -and $h.Range.Column -eq $Col) {
            $h.TextToDisplay = $Value
        }
    }
}

# ---------------------------------------------------------------------------
# Sheet "Overview": columns A=File Name, B=zh-cn, C=de-de, D=Latest Handoff Date
# ---------------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

Set-CellAndLink $wsOverview 2 1 "ffff6c0ace36-61ad-4947-a5d9-19faa014a1d7.md"
Set-CellAndLink $wsOverview 2 2 "Handed back: in sync with en-US"
Set-CellAndLink $wsOverview 2 3 "Handed back: in sync with en-US"
Set-CellAndLink $wsOverview 2 4 "2016-03-22 21:11:02"

Set-CellAndLink $wsOverview 3 1 "ffffffab0b8d44-08cb-4112-bb04-99628528bfad.md"
Set-CellAndLink $wsOverview 3 2 "Handed back: in sync with en-US"
Set-CellAndLink $wsOverview 3 3 "Handed back: in sync with en-US"
Set-CellAndLink $wsOverview 3 4 "2016-03-22 21:11:02"

Set-CellAndLink $wsOverview 4 1 "faec0f94-5a0c-4cd5-a6b8-19c4b111ec72.md"
Set-CellAndLink $wsOverview 4 2 "Ready for handoff"
Set-CellAndLink $wsOverview 4 3 "Ready for handoff"
Set-CellAndLink $wsOverview 4 4 "2016-03-22 21:13:55"

# ---------------------------------------------------------------------------
# Sheet "zh-cn": A=Source File Name, B=File Extension, C=Status,
# D=Latest Handoff File, E=Latest Handoff Datetime, F=Latest Target File,
# G=Latest Handback File, H=Latest Handback DateTime, J=Handoff Reason
# ---------------------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")

Set-CellAndLink $wsZh 2 1 "ffff6c0ace36-61ad-4947-a5d9-19faa014a1d7.md"
Set-CellAndLink $wsZh 2 2 ".md"
Set-CellAndLink $wsZh 2 3 "Handed back: in sync with en-US"
Set-CellAndLink $wsZh 2 4 "c85dd3f2-ec34-4d7e-980f-d33a2a27cfba.eb7edf2a86468b4bd614ade89c8221dda5c35aab.zh-cn.xlf"
Set-CellAndLink $wsZh 2 5 "2016-03-22 21:10:58"
Set-CellAndLink $wsZh 2 6 "c85dd3f2-ec34-4d7e-980f-d33a2a27cfba.md"
Set-CellAndLink $wsZh 2 7 "c85dd3f2-ec34-4d7e-980f-d33a2a27cfba.eb7edf2a86468b4bd614ade89c8221dda5c35aab.zh-cn.xlf"
Set-CellAndLink $wsZh 2 8 "2016-03-22 21:11:29"
Set-CellAndLink $wsZh 2 10 "Include"

Set-CellAndLink $wsZh 3 1 "ffffffab0b8d44-08cb-4112-bb04-99628528bfad.md"
Set-CellAndLink $wsZh 3 2 ".md"
Set-CellAndLink $wsZh 3 3 "Handed back: in sync with en-US"
Set-CellAndLink $wsZh 3 4 "c85dd3f2-ec34-4d7e-980f-d33a2a27cfba.eb7edf2a86468b4bd614ade89c8221dda5c35aab.zh-cn.xlf"
Set-CellAndLink $wsZh 3 5 "2016-03-22 21:10:58"
Set-CellAndLink $wsZh 3 6 "c85dd3f2-ec34-4d7e-980f-d33a2a27cfba.md"
Set-CellAndLink $wsZh 3 7 "c85dd3f2-ec34-4d7e-980f-d33a2a27cfba.eb7edf2a86468b4bd614ade89c8221dda5c35aab.zh-cn.xlf"
Set-CellAndLink $wsZh 3 8 "2016-03-22 21:11:29"
Set-CellAndLink $wsZh 3 10 "Include"

Set-CellAndLink $wsZh 4 1 "faec0f94-5a0c-4cd5-a6b8-19c4b111ec72.md"
Set-CellAndLink $wsZh 4 2 ".md"
Set-CellAndLink $wsZh 4 3 "Ready for handoff"
Set-CellAndLink $wsZh 4 4 "faec0f94-5a0c-4cd5-a6b8-19c4b111ec72.d78a5ad2165d71dde2790e69cf593f42785c00ee.zh-cn.xlf"
Set-CellAndLink $wsZh 4 5 "2016-03-22 21:13:51"
Set-CellAndLink $wsZh 4 6 "faec0f94-5a0c-4cd5-a6b8-19c4b111ec72.md"
Set-CellAndLink $wsZh 4 7 "faec0f94-5a0c-4cd5-a6b8-19c4b111ec72.d78a5ad2165d71dde2790e69cf593f42785c00ee.zh-cn.xlf"
Set-CellAndLink $wsZh 4 8 "2016-03-22 21:13:04"
Set-CellAndLink $wsZh 4 10 "Include"

# ---------------------------------------------------------------------------
# Sheet "de-de": same column layout as "zh-cn"
# ---------------------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")

Set-CellAndLink $wsDe 2 1 "ffff6c0ace36-61ad-4947-a5d9-19faa014a1d7.md"
Set-CellAndLink $wsDe 2 2 ".md"
Set-CellAndLink $wsDe 2 3 "Handed back: in sync with en-US"
Set-CellAndLink $wsDe 2 4 "c85dd3f2-ec34-4d7e-980f-d33a2a27cfba.eb7edf2a86468b4bd614ade89c8221dda5c35aab.de-de.xlf"
Set-CellAndLink $wsDe 2 5 "2016-03-22 21:11:02"
Set-CellAndLink $wsDe 2 6 "c85dd3f2-ec34-4d7e-980f-d33a2a27cfba.md"
Set-CellAndLink $wsDe 2 7 "c85dd3f2-ec34-4d7e-980f-d33a2a27cfba.eb7edf2a86468b4bd614ade89c8221dda5c35aab.de-de.xlf"
Set-CellAndLink $wsDe 2 8 "2016-03-22 21:11:38"
Set-CellAndLink $wsDe 2 10 "Include"

Set-CellAndLink $wsDe 3 1 "ffffffab0b8d44-08cb-4112-bb04-99628528bfad.md"
Set-CellAndLink $wsDe 3 2 ".md"
Set-CellAndLink $wsDe 3 3 "Handed back: in sync with en-US"
Set-CellAndLink $wsDe 3 4 "c85dd3f2-ec34-4d7e-980f-d33a2a27cfba.eb7edf2a86468b4bd614ade89c8221dda5c35aab.de-de.xlf"
Set-CellAndLink $wsDe 3 5 "2016-03-22 21:11:02"
Set-CellAndLink $wsDe 3 6 "c85dd3f2-ec34-4d7e-980f-d33a2a27cfba.md"
Set-CellAndLink $wsDe 3 7 "c85dd3f2-ec34-4d7e-980f-d33a2a27cfba.eb7edf2a86468b4bd614ade89c8221dda5c35aab.de-de.xlf"
Set-CellAndLink $wsDe 3 8 "2016-03-22 21:11:38"
Set-CellAndLink $wsDe 3 10 "Include"

Set-CellAndLink $wsDe 4 1 "faec0f94-5a0c-4cd5-a6b8-19c4b111ec72.md"
Set-CellAndLink $wsDe 4 2 ".md"
Set-CellAndLink $wsDe 4 3 "Ready for handoff"
Set-CellAndLink $wsDe 4 4 "faec0f94-5a0c-4cd5-a6b8-19c4b111ec72.d78a5ad2165d71dde2790e69cf593f42785c00ee.de-de.xlf"
Set-CellAndLink $wsDe 4 5 "2016-03-22 21:13:55"
Set-CellAndLink $wsDe 4 6 "faec0f94-5a0c-4cd5-a6b8-19c4b111ec72.md"
Set-CellAndLink $wsDe 4 7 "faec0f94-5a0c-4cd5-a6b8-19c4b111ec72.d78a5ad2165d71dde2790e69cf593f42785c00ee.de-de.xlf"
Set-CellAndLink $wsDe 4 8 "2016-03-22 21:13:10"
Set-CellAndLink $wsDe 4 10 "Include"
